# Generate Report for Handoff
# Adds a new handed-off file (fa34a3a6-f750-4c1b-8679-9b2bd082f09a.md) as a new
# row (row 3) on the Overview, zh-cn and de-de sheets, mirroring the existing
# c54ca49b-... row, and resizes the worksheet tables to include the new row.

$wb = $excel.ActiveWorkbook

$commit = "f65e56b2f7ea6dff79adab564f6d9b1d2af9af98"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/"

$newFile = "fa34a3a6-f750-4c1b-8679-9b2bd082f09a.md"
$newFileDisplayPath = "e2e\$newFile"
$newZhXlf = "fa34a3a6-f750-4c1b-8679-9b2bd082f09a.694cb695c79c3471acbbfd66730ea1b312941fae.zh-cn.xlf"
$newDeXlf = "fa34a3a6-f750-4c1b-8679-9b2bd082f09a.694cb695c79c3471acbbfd66730ea1b312941fae.de-de.xlf"

$zhHandoffDate = "2016-08-18 06:38:28"
$deHandoffDate = "2016-08-18 06:38:33"
$overviewDate = "2016-08-18 06:38:33"

$dateFormat = "yyyy-mm-dd HH:mm:ss"
# OLE color equivalent of the workbook's existing HyperLink font color (FF6495ED)
$hyperlinkColor = 15570276

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newFileDisplayPath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($baseUrl + $newFile), "", "", $newFileDisplayPath)
Set-HyperlinkLook($wsOverview.Range("B3"))
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = $dateFormat

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - new row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newFile
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($baseUrl + $newFile), "", "", $newFile)
Set-HyperlinkLook($wsZh.Range("A3"))
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("H3").NumberFormat = $dateFormat
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $dateFormat
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - new row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newFile
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($baseUrl + $newFile), "", "", $newFile)
Set-HyperlinkLook($wsDe.Range("A3"))
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("H3").NumberFormat = $dateFormat
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $dateFormat
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Host "Handoff report row added for $newFile"
